$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The test data list used to contain 4 products (Backpack, Bike Light,
# Bolt T-Shirt, Onesie). The alert-handling test only needs a single
# product now, so trim the list down to just "Sauce Labs Onesie",
# immediately below the "Name" header.
$ws.Range("A2").Value = "Sauce Labs Onesie"

# Delete the now-unneeded rows (old rows 3-5: Bike Light, Bolt T-Shirt,
# and the old Onesie row) so the remaining data shifts up and the used
# range shrinks to A1:A2.
$ws.Rows("3:5").Delete() | Out-Null

# Excel leaves the active selection right below the last data row after
# such edits.
$ws.Range("A2").Select() | Out-Null
